# Refresh this NATMI LR-pair sheet with recomputed TPM-based values.
# The underlying TPM table changed (notably ligand expression for the
# "ECs" sending cluster and receptor-expressing-cell counts/levels for the
# "ECs" target cluster), so every derived column (G..T) is refreshed below
# with the recomputed values, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs
$ws.Cells.Item(2, 7).Value = 0.1887043333333333  # G - Ligand average expression value
$ws.Cells.Item(2, 8).Value = 0.566113  # H - Ligand total expression value
$ws.Cells.Item(2, 9).Value = 0.02109097403787168  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(2, 10).Value = 0.02109097403787168  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(2, 11).Value = 3  # K - Receptor-expressing cells
$ws.Cells.Item(2, 12).Value = 1  # L - Receptor detection rate
$ws.Cells.Item(2, 13).Value = 28.22405966666667  # M - Receptor average expression value
$ws.Cells.Item(2, 14).Value = 84.672179  # N - Receptor total expression value
$ws.Cells.Item(2, 15).Value = 0.3816548478108986  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(2, 16).Value = 0.3816548478108986  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(2, 17).Value = 5.326002363358556  # Q - Edge average expression weight
$ws.Cells.Item(2, 18).Value = 47.934021270227  # R - Edge total expression weight
$ws.Cells.Item(2, 19).Value = 0.008049472486607532  # S - Edge average expression derived specificity
$ws.Cells.Item(2, 20).Value = 0.008049472486607529  # T - Edge total expression derived specificity

# Row 3: Sending=ECs, Target=FAPs
$ws.Cells.Item(3, 7).Value = 0.1887043333333333  # G - Ligand average expression value
$ws.Cells.Item(3, 8).Value = 0.566113  # H - Ligand total expression value
$ws.Cells.Item(3, 9).Value = 0.02109097403787168  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(3, 10).Value = 0.02109097403787168  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(3, 14).Value = 59.306181  # N - Receptor total expression value
$ws.Cells.Item(3, 15).Value = 0.2673191094302723  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(3, 16).Value = 0.2673191094302723  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(3, 17).Value = 3.730444449383667  # Q - Edge average expression weight
$ws.Cells.Item(3, 18).Value = 33.57400004445299  # R - Edge total expression weight
$ws.Cells.Item(3, 19).Value = 0.005638020396820853  # S - Edge average expression derived specificity
$ws.Cells.Item(3, 20).Value = 0.005638020396820851  # T - Edge total expression derived specificity

# Row 4: Sending=ECs, Target=MuSCs
$ws.Cells.Item(4, 7).Value = 0.1887043333333333  # G - Ligand average expression value
$ws.Cells.Item(4, 8).Value = 0.566113  # H - Ligand total expression value
$ws.Cells.Item(4, 9).Value = 0.02109097403787168  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(4, 10).Value = 0.02109097403787168  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(4, 13).Value = 25.95900466666667  # M - Receptor average expression value
$ws.Cells.Item(4, 14).Value = 77.877014  # N - Receptor total expression value
$ws.Cells.Item(4, 15).Value = 0.351026042758829  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(4, 16).Value = 0.351026042758829  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(4, 17).Value = 4.898576669620223  # Q - Edge average expression weight
$ws.Cells.Item(4, 18).Value = 44.087190026582  # R - Edge total expression weight
$ws.Cells.Item(4, 19).Value = 0.007403481154443298  # S - Edge average expression derived specificity
$ws.Cells.Item(4, 20).Value = 0.007403481154443295  # T - Edge total expression derived specificity

# Row 5: Sending=FAPs, Target=ECs
$ws.Cells.Item(5, 7).Value = 7.8617  # G - Ligand average expression value
$ws.Cells.Item(5, 9).Value = 0.8786809908633213  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(5, 10).Value = 0.8786809908633211  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(5, 11).Value = 3  # K - Receptor-expressing cells
$ws.Cells.Item(5, 12).Value = 1  # L - Receptor detection rate
$ws.Cells.Item(5, 13).Value = 28.22405966666667  # M - Receptor average expression value
$ws.Cells.Item(5, 14).Value = 84.672179  # N - Receptor total expression value
$ws.Cells.Item(5, 15).Value = 0.3816548478108986  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(5, 16).Value = 0.3816548478108986  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(5, 17).Value = 221.8890898814333  # Q - Edge average expression weight
$ws.Cells.Item(5, 18).Value = 1997.0018089329  # R - Edge total expression weight
$ws.Cells.Item(5, 19).Value = 0.3353528598422705  # S - Edge average expression derived specificity
$ws.Cells.Item(5, 20).Value = 0.3353528598422704  # T - Edge total expression derived specificity

# Row 6: Sending=FAPs, Target=FAPs
$ws.Cells.Item(6, 7).Value = 7.8617  # G - Ligand average expression value
$ws.Cells.Item(6, 9).Value = 0.8786809908633213  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(6, 10).Value = 0.8786809908633211  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(6, 14).Value = 59.306181  # N - Receptor total expression value
$ws.Cells.Item(6, 15).Value = 0.2673191094302723  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(6, 16).Value = 0.2673191094302723  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(6, 19).Value = 0.2348882199508923  # S - Edge average expression derived specificity
$ws.Cells.Item(6, 20).Value = 0.2348882199508922  # T - Edge total expression derived specificity

# Row 7: Sending=FAPs, Target=MuSCs
$ws.Cells.Item(7, 7).Value = 7.8617  # G - Ligand average expression value
$ws.Cells.Item(7, 9).Value = 0.8786809908633213  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(7, 10).Value = 0.8786809908633211  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(7, 13).Value = 25.95900466666667  # M - Receptor average expression value
$ws.Cells.Item(7, 14).Value = 77.877014  # N - Receptor total expression value
$ws.Cells.Item(7, 15).Value = 0.351026042758829  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(7, 16).Value = 0.351026042758829  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(7, 19).Value = 0.3084399110701585  # S - Edge average expression derived specificity
$ws.Cells.Item(7, 20).Value = 0.3084399110701584  # T - Edge total expression derived specificity

# Row 8: Sending=MuSCs, Target=ECs
$ws.Cells.Item(8, 9).Value = 0.1002280350988072  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(8, 10).Value = 0.1002280350988071  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(8, 11).Value = 3  # K - Receptor-expressing cells
$ws.Cells.Item(8, 12).Value = 1  # L - Receptor detection rate
$ws.Cells.Item(8, 13).Value = 28.22405966666667  # M - Receptor average expression value
$ws.Cells.Item(8, 14).Value = 84.672179  # N - Receptor total expression value
$ws.Cells.Item(8, 15).Value = 0.3816548478108986  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(8, 16).Value = 0.3816548478108986  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(8, 17).Value = 25.31010425846122  # Q - Edge average expression weight
$ws.Cells.Item(8, 18).Value = 227.790938326151  # R - Edge total expression weight
$ws.Cells.Item(8, 19).Value = 0.03825251548202065  # S - Edge average expression derived specificity
$ws.Cells.Item(8, 20).Value = 0.03825251548202064  # T - Edge total expression derived specificity

# Row 9: Sending=MuSCs, Target=FAPs
$ws.Cells.Item(9, 9).Value = 0.1002280350988072  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(9, 10).Value = 0.1002280350988071  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(9, 14).Value = 59.306181  # N - Receptor total expression value
$ws.Cells.Item(9, 15).Value = 0.2673191094302723  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(9, 16).Value = 0.2673191094302723  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(9, 17).Value = 17.72773113918766  # Q - Edge average expression weight
$ws.Cells.Item(9, 19).Value = 0.0267928690825592  # S - Edge average expression derived specificity
$ws.Cells.Item(9, 20).Value = 0.0267928690825592  # T - Edge total expression derived specificity

# Row 10: Sending=MuSCs, Target=MuSCs
$ws.Cells.Item(10, 9).Value = 0.1002280350988072  # I - Ligand derived specificity of average expression value
$ws.Cells.Item(10, 10).Value = 0.1002280350988071  # J - Ligand derived specificity of total expression value
$ws.Cells.Item(10, 13).Value = 25.95900466666667  # M - Receptor average expression value
$ws.Cells.Item(10, 14).Value = 77.877014  # N - Receptor total expression value
$ws.Cells.Item(10, 15).Value = 0.351026042758829  # O - Receptor derived specificity of average expression value
$ws.Cells.Item(10, 16).Value = 0.351026042758829  # P - Receptor derived specificity of total expression value
$ws.Cells.Item(10, 17).Value = 23.27890184186289  # Q - Edge average expression weight
$ws.Cells.Item(10, 19).Value = 0.0351826505342273  # S - Edge average expression derived specificity
$ws.Cells.Item(10, 20).Value = 0.0351826505342273  # T - Edge total expression derived specificity
